$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2
Set-TextCell "D2" "42.024.33"
Set-TextCell "E2" "  -3.83%  "

# Row 3
Set-TextCell "D3" "2.194.32"
Set-TextCell "E3" "  -3.77%  "

# Row 4
Set-TextCell "E4" "  +0.53%  "

# Row 5
Set-TextCell "D5" "106.35"
Set-TextCell "E5" "  -14.49%  "

# Row 6
Set-TextCell "D6" "290.89"
Set-TextCell "E6" "  +9.34%  "

# Row 7
Set-TextCell "D7" "0.616"
Set-TextCell "E7" "  -3.50%  "

# Row 8
Set-TextCell "D8" "1.00"
Set-TextCell "E8" "  -0.02%  "

# Row 9
Set-TextCell "E9" "  -6.48%  "

# Row 10
Set-TextCell "D10" "43.44"
Set-TextCell "E10" "  -9.77%  "

# Row 11
Set-TextCell "D11" "0.0901"
Set-TextCell "E11" "  -4.75%  "

# Row 12
Set-TextCell "D12" "54.01"
Set-TextCell "E12" "  -0.55%  "

# Row 13
Set-TextCell "D13" "8.61"
Set-TextCell "E13" "  -6.82%  "

# Row 14
Set-TextCell "E14" "  -3.54%  "

# Row 15
Set-TextCell "D15" "0.924"
Set-TextCell "E15" "  +2.73%  "

# Row 16
Set-TextCell "D16" "14.70"
Set-TextCell "E16" "  -4.78%  "

# Row 17
Set-TextCell "D17" "2.531.96"
Set-TextCell "E17" "  -3.49%  "

# Row 18
Set-TextCell "D18" "2.229.30"
Set-TextCell "E18" "  -2.12%  "

# Row 19
Set-TextCell "D19" "41.999.48"
Set-TextCell "E19" "  -3.85%  "

# Row 20
Set-TextCell "D20" "7.11"
Set-TextCell "E20" "  +1.35%  "

# Row 21
Set-TextCell "E21" "  -6.19%  "

# Row 22
Set-TextCell "D22" "72.06"
Set-TextCell "E22" "  -0.47%  "

# Row 23
Set-TextCell "D23" "3.36"
Set-TextCell "E23" "  +16.54%  "

# Row 24
Set-TextCell "D24" "2.24"
Set-TextCell "E24" "  -8.23%  "

# Row 25
Set-TextCell "D25" "225.24"
Set-TextCell "E25" "  -4.43%  "

# Row 26
Set-TextCell "D26" "8.86"
Set-TextCell "E26" "  -6.17%  "

# Row 27
Set-TextCell "D27" "0.999"
Set-TextCell "E27" "  -1.78%  "

# Row 28
Set-TextCell "D28" "11.42"
Set-TextCell "E28" "  -3.45%  "

# Row 29
Set-TextCell "D29" "3.89"
Set-TextCell "E29" "  -1.11%  "

# Row 30
Set-TextCell "E30" "  -1.75%  "

# Row 31
Set-TextCell "E31" "  -5.11%  "

# Row 32
Set-TextCell "D32" "37.07"
Set-TextCell "E32" "  -13.31%  "

# Row 33
Set-TextCell "D33" "170.87"
Set-TextCell "E33" "  -1.09%  "

# Row 34
Set-TextCell "D34" "20.58"
Set-TextCell "E34" "  -5.12%  "

# Row 35
Set-TextCell "D35" "0.0857"
Set-TextCell "E35" "  -6.40%  "

# Row 36
Set-TextCell "D36" "5.44"
Set-TextCell "E36" "  -5.98%  "

# Row 37
Set-TextCell "D37" "4.81"
Set-TextCell "E37" "  +3.26%  "

# Row 38
Set-TextCell "D38" "4.15"
Set-TextCell "E38" "  -1.70%  "

# Row 39
Set-TextCell "E39" "  -4.46%  "

# Row 40
Set-TextCell "D40" "0.0356"
Set-TextCell "E40" "  -5.62%  "

# Row 41
Set-TextCell "E41" "  -6.53%  "

# Row 42
Set-TextCell "D42" "2.40"
Set-TextCell "E42" "  -5.55%  "

# Row 43
$ws.Range("B43").Value = "MultiversX"
$ws.Range("C43").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
Set-TextCell "D43" "68.96"
Set-TextCell "E43" "  -8.23%  "

# Row 44
$ws.Range("B44").Value = "Algorand"
$ws.Range("C44").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextCell "D44" "0.226"
Set-TextCell "E44" "  -5.67%  "

# Row 45
Set-TextCell "D45" "1.00"
Set-TextCell "E45" "  +0.23%  "

# Row 46
Set-TextCell "D46" "12.40"
Set-TextCell "E46" "  -11.20%  "

# Row 47
Set-TextCell "E47" "  -7.04%  "

# Row 48
Set-TextCell "D48" "5.33"
Set-TextCell "E48" "  -4.41%  "

# Row 49
Set-TextCell "D49" "1.28"
Set-TextCell "E49" "  +1.17%  "

# Row 50
Set-TextCell "D50" "101.04"
Set-TextCell "E50" "  -0.85%  "

# Row 51
Set-TextCell "D51" "8.26"
Set-TextCell "E51" "  -4.08%  "
